# Generate Report for Handoff
# The "b.md" row moves from "Handed back: in sync with en-US" to
# "Ready for handoff": a new handoff xliff was generated for b.md, and the
# previous handback for b.md is now stale (content-duplicate check failed,
# producing an "Error Detail" message).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is b.md ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Cells.Item(3, 5).Value = "Ready for handoff"              # E3 (zh-cn status)
$ov.Cells.Item(3, 6).Value = "Ready for handoff"              # F3 (de-de status)
$ov.Cells.Item(3, 7).Value = "2016-08-30 00:39:26"             # G3 (latest HO xliff generate date)

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09fb7839915c6404d4984363e6cf3b07567038a1/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5431f1353ee906e5b310dcd056fcc9d9554d9aa/e2e/b.md."

# ---- zh-cn sheet: row 3 is b.md ----
# NOTE: "'False" (leading apostrophe) forces a text cell instead of letting
# Excel auto-coerce the literal word False/True into a Boolean cell type;
# re-applying the "Normal" style afterwards drops the quote-prefix styling
# so the cell ends up plain text with the default (unstyled) look, same as
# its neighbours.
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Cells.Item(3, 3).Value = "Ready for handoff"                                               # C3 Status
$zh.Cells.Item(3, 6).Value = "'False"                                                          # F3 Content Duplicate
$zh.Cells.Item(3, 6).Style = "Normal"
$zh.Cells.Item(3, 7).Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"             # G3 Latest Handoff File
$zh.Cells.Item(3, 8).Value = "2016-08-30 00:39:21"                                             # H3 Latest Handoff Datetime
$zh.Cells.Item(3, 16).Value = $errorDetail                                                     # P3 Error Detail
$zh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 is b.md ----
$de = $wb.Worksheets.Item("de-de")
$de.Cells.Item(3, 3).Value = "Ready for handoff"                                               # C3 Status
$de.Cells.Item(3, 6).Value = "'False"                                                          # F3 Content Duplicate
$de.Cells.Item(3, 6).Style = "Normal"
$de.Cells.Item(3, 7).Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"             # G3 Latest Handoff File
$de.Cells.Item(3, 8).Value = "2016-08-30 00:39:26"                                             # H3 Latest Handoff Datetime
$de.Cells.Item(3, 16).Value = $errorDetail                                                     # P3 Error Detail
$de.Columns.Item(16).ColumnWidth = 39.17
